$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): insert new "id" column at A, shift others, drop old "margem_venda",
#     add new "preco_venda" earlier and append "preco_custo" at the end ---
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "codigo_barras"
$ws.Range("C1").Value = "nome"
$ws.Range("D1").Value = "tipo"
$ws.Range("E1").Value = "preco_venda"
$ws.Range("F1").Value = "estoque"
$ws.Range("G1").Value = "categoria"
$ws.Range("H1").Value = "perecivel"
$ws.Range("I1").Value = "validade"
$ws.Range("J1").Value = "preco_custo"

# --- Data row (row 2): now holds the "Alface Und" product, remapped to the new columns ---
$ws.Range("A2").Value = 1
# codigo_barras "100" must stay textual (not coerce to a number)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "100"
$ws.Range("B2").ClearFormats()
$ws.Range("C2").Value = "Alface Und"
$ws.Range("D2").Value = "unidade"
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 14
$ws.Range("G2").Value = "Hortaliças"
$ws.Range("H2").Value = "nao"
# I2 / J2 stay untouched - they were (and remain) empty cells

# --- Drop the old column K (id used to live there) and the now-unused rows 3-6 ---
$ws.Columns("K").Delete()
$ws.Rows("3:6").Delete()
